$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "groteauto"
$ws.Range("C2").Value = "GROTE AUTOMATION"

$ws.Range("C3").Select()
